# Update "Base de données" exercise table: add field lengths to the
# "Type (Longueur)" column and fix a couple of mnemonic typos/renames.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Row 2: beach_id — Entier -> Numérique (10) --------------------------
$t.Cell(2, 3).Range.Text = "Numérique (10)"

# --- Row 3: beach_name — Alphabétique -> Alphabétique (50) ---------------
$t.Cell(3, 3).Range.Text = "Alphabétique (50)"

# --- Row 4: beach_lenght -> beach_length (typo fix) -----------------------
# The mnemonic cell has several runs ("b" / "each" / "_l" / "enght") plus
# proofing-error markers, so only touch the final run's text ("enght" ->
# "ength") instead of overwriting the whole multi-run range.
$cell = $t.Cell(4, 1)
$cellStart = $cell.Range.Start
$sub = $d.Range($cellStart + 7, $cellStart + 12)   # "enght"
$sub.Text = "ength"

# --- Row 4: Décimal -> Numérique (10,2) -----------------------------------
$t.Cell(4, 3).Range.Text = "Numérique (10,2)"

# --- Row 5: type_of_land — Alphabétique -> Alphabétique (15) -------------
$t.Cell(5, 3).Range.Text = "Alphabétique (15)"

# --- Row 6: city — Alphabétique -> Alphabétique (20) ----------------------
$t.Cell(6, 3).Range.Text = "Alphabétique (20)"

# --- Row 7: postal_code — Entier -> Numérique (5) -------------------------
$t.Cell(7, 3).Range.Text = "Numérique (5)"

# --- Row 8: regional_department — Alphabétique -> Alphabétique (20) ------
$t.Cell(8, 3).Range.Text = "Alphabétique (20)"

# --- Row 9: manager_name -> manager_lastname ------------------------------
$cell = $t.Cell(9, 1)
$cellStart = $cell.Range.Start
$sub = $d.Range($cellStart + 7, $cellStart + 9)    # "_n"
$sub.Text = "_lastn"
$t.Cell(9, 3).Range.Text = "Alphabétique (50)"

# --- Row 10: manager_first_name -> manager_firstname ----------------------
$cell = $t.Cell(10, 1)
$cellStart = $cell.Range.Start
$sub = $d.Range($cellStart + 13, $cellStart + 15)  # "_n"
$sub.Text = "n"
$t.Cell(10, 3).Range.Text = "Alphabétique (50)"

# --- Row 11: annual_nb_of_tourist — Entier -> Numérique (10) -------------
$t.Cell(11, 3).Range.Text = "Numérique (10)"
